$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Logs")
$dash = $wb.Worksheets.Item("Dashboard")

# Add the new log row (row 3)
$ws.Range("A3").Value = "Retour status"
$ws.Range("B3").Value = "mailmind.test@zohomail.eu"
$ws.Range("D3").Value = "Klantenservice / Opvolging"
$ws.Range("F3").Value = "2025-08-22 22:34:01"
$ws.Range("G3").Value = "Nee"
$ws.Range("H3").Value = "Ja"
$ws.Range("I3").Value = "Nee"
$ws.Range("J3").Value = "Nee"

# Extend the conditional formatting ranges to cover the new row
$ws.Range("D2").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("D2:D3"))
$ws.Range("G2").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("G2:G3"))
$ws.Range("H2").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("H2:H3"))
$ws.Range("I2").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("I2:I3"))
$ws.Range("J2").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("J2:J3"))

# Update the dashboard count for this category
$dash.Range("B2").Value = 2
